# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Appends 6 new transaction rows (89-94) to Sheet1, matching the new
# shared-string entries for Date, Number and Name columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 89; Date = "2026-02-19 12:16:35"; Number = "237674890585"; Name = "LA NEGRESSE LTDLA CBOX R0 NGO MBOCK epse MBAYAN MARIE CLAIRE"; Balance = 67456 },
    @{ Row = 90; Date = "2026-02-20 10:49:58"; Number = "237653294562"; Name = "NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS"; Balance = 506519 },
    @{ Row = 91; Date = "2026-02-20 15:15:43"; Number = "237681659043"; Name = "SYLVIE DJIDJOU TEGUIA EPSE TOUKOU"; Balance = 18480 },
    @{ Row = 92; Date = "2026-02-20 16:42:30"; Number = "237679086144"; Name = "ALAIN CHETEU KAMDEM"; Balance = 2615 },
    @{ Row = 93; Date = "2026-02-20 11:17:54"; Number = "237683998069"; Name = "MEDJEU FEUZEU epse FEGHEM WAHOUE JOSIANE ETS MOBILE FINANCIAL SERVICES MFS"; Balance = 17949 },
    @{ Row = 94; Date = "2026-02-17 12:41:34"; Number = "237679252522"; Name = "WARAMMA NICOLAS"; Balance = 23 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Column A - Date (kept as plain text, matches existing column formatting)
    $ws.Cells.Item($rowIndex, 1).Value = $r.Date

    # Column B - Number (12-digit string that must stay text, not be
    # coerced to a numeric value). Force text format, assign, then
    # restore the default "Normal" style so no stray number format is
    # left attached to the cell.
    $cellB = $ws.Cells.Item($rowIndex, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.Number
    $cellB.Style = "Normal"

    # Column C - Name (plain text)
    $ws.Cells.Item($rowIndex, 3).Value = $r.Name

    # Column D - Balance (numeric)
    $ws.Cells.Item($rowIndex, 4).Value = $r.Balance
}
